$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("C2").Value = 4.926126621258732
$ws.Range("D2").Value = 3.338002904462885
$ws.Range("E2").Value = 16.66695836158539
$ws.Range("F2").Value = 27.18703899061369
$ws.Range("G2").Value = 3.567828460017615
$ws.Range("N2").Value = 18.99769123883711
$ws.Range("O2").Value = 23.19951237517529

# Row 3
$ws.Range("C3").Value = 4.750213261576305
$ws.Range("D3").Value = 3.332266480516493
$ws.Range("E3").Value = 15.70294774289624
$ws.Range("F3").Value = 26.3124813607043
$ws.Range("G3").Value = 3.572303459151253
$ws.Range("N3").Value = 18.40031508502701
$ws.Range("O3").Value = 22.57619591036837

# Row 4
$ws.Range("C4").Value = 4.640614124970813
$ws.Range("D4").Value = 3.329717909951302
$ws.Range("E4").Value = 15.08611945890486
$ws.Range("F4").Value = 25.77068349195066
$ws.Range("G4").Value = 3.575189172148977
$ws.Range("N4").Value = 18.02485520896362
$ws.Range("O4").Value = 22.19320626190153

# Row 5
$ws.Range("C5").Value = 4.595631543631624
$ws.Range("D5").Value = 3.328923112997993
$ws.Range("E5").Value = 14.82875600235973
$ws.Range("F5").Value = 25.54907308206959
$ws.Range("G5").Value = 3.576399996207454
$ws.Range("N5").Value = 17.86990355188765
$ws.Range("O5").Value = 22.03732681160587

# Row 6
$ws.Range("C6").Value = 4.58814532945987
$ws.Range("D6").Value = 3.328805825205932
$ws.Range("E6").Value = 14.78566738895239
$ws.Range("F6").Value = 25.51223691030438
$ws.Range("G6").Value = 3.576603163521045
$ws.Range("N6").Value = 17.84406337566597
$ws.Range("O6").Value = 22.01146251924026

# Row 7
$ws.Range("C7").Value = 4.640008660651212
$ws.Range("D7").Value = 3.329706205594031
$ws.Range("E7").Value = 15.08267247213796
$ws.Range("F7").Value = 25.76769758027975
$ws.Range("G7").Value = 3.575205360334442
$ws.Range("N7").Value = 18.02277304766462
$ws.Range("O7").Value = 22.19110288947427

# Row 8
$ws.Range("C8").Value = 4.865846705030525
$ws.Range("D8").Value = 3.33582233936324
$ws.Range("E8").Value = 16.33988363937716
$ws.Range("F8").Value = 26.88671595412645
$ws.Range("G8").Value = 3.569342884976606
$ws.Range("N8").Value = 18.79364780656866
$ws.Range("O8").Value = 22.98480201195844

# Row 9
$ws.Range("C9").Value = 5.292972261085342
$ws.Range("D9").Value = 3.355579124847327
$ws.Range("E9").Value = 18.74361448004077
$ws.Range("F9").Value = 29.02717102200151
$ws.Range("G9").Value = 3.55893461158872
$ws.Range("N9").Value = 20.2273683202997
$ws.Range("O9").Value = 24.52821902330702

# Row 10
$ws.Range("C10").Value = 5.593444883520045
$ws.Range("D10").Value = 3.374865229620801
$ws.Range("E10").Value = 20.43372571100995
$ws.Range("F10").Value = 30.54778217139864
$ws.Range("G10").Value = 3.551940698096631
$ws.Range("N10").Value = 21.22223697909767
$ws.Range("O10").Value = 25.64059188266682

# Row 11
$ws.Range("C11").Value = 5.72658274896879
$ws.Range("D11").Value = 3.384678399948793
$ws.Range("E11").Value = 21.16060891215216
$ws.Range("F11").Value = 31.22479402351992
$ws.Range("G11").Value = 3.54889858406329
$ws.Range("N11").Value = 21.66018057919901
$ws.Range("O11").Value = 26.13936855543931

# Row 12
$ws.Range("C12").Value = 5.776441726129652
$ws.Range("D12").Value = 3.388543798539379
$ws.Range("E12").Value = 21.42987029500256
$ws.Range("F12").Value = 31.47879319919046
$ws.Range("G12").Value = 3.547766495598579
$ws.Range("N12").Value = 21.82377585682186
$ws.Range("O12").Value = 26.32700660914469

# Row 13
$ws.Range("C13").Value = 5.765729185756575
$ws.Range("D13").Value = 3.387704677737943
$ws.Range("E13").Value = 21.37214596823456
$ws.Range("F13").Value = 31.42419911024609
$ws.Range("G13").Value = 3.548009428843297
$ws.Range("N13").Value = 21.78864458690801
$ws.Range("O13").Value = 26.28665339029789

# Row 14
$ws.Range("C14").Value = 5.730696098716328
$ws.Range("D14").Value = 3.384993416870911
$ws.Range("E14").Value = 21.18288109909063
$ws.Range("F14").Value = 31.2457395530387
$ws.Range("G14").Value = 3.548805048668281
$ws.Range("N14").Value = 21.67368539489659
$ws.Range("O14").Value = 26.15483152136503

# Row 15
$ws.Range("C15").Value = 5.709163424527925
$ws.Range("D15").Value = 3.383352134246597
$ws.Range("E15").Value = 21.06617155094796
$ws.Range("F15").Value = 31.13611214935311
$ws.Range("G15").Value = 3.54929497485921
$ws.Range("N15").Value = 21.60297336126124
$ws.Range("O15").Value = 26.07392009104949

# Row 16
$ws.Range("C16").Value = 5.584668440015143
$ws.Range("D16").Value = 3.374244831144299
$ws.Range("E16").Value = 20.38538131221762
$ws.Range("F16").Value = 30.50321931614264
$ws.Range("G16").Value = 3.552142299956858
$ws.Range("N16").Value = 21.19330956972086
$ws.Range("O16").Value = 25.60783271292764

# Row 17
$ws.Range("C17").Value = 5.507350503590494
$ws.Range("D17").Value = 3.368924076746687
$ws.Range("E17").Value = 19.95702076391906
$ws.Range("F17").Value = 30.11100033487716
$ws.Range("G17").Value = 3.553924645731941
$ws.Range("N17").Value = 20.93814219015166
$ws.Range("O17").Value = 25.31990249809792

# Row 18
$ws.Range("C18").Value = 5.462547692818903
$ws.Range("D18").Value = 3.36596161706636
$ws.Range("E18").Value = 19.70668754964532
$ws.Range("F18").Value = 29.88403404494041
$ws.Range("G18").Value = 3.554962939315403
$ws.Range("N18").Value = 20.79000725568362
$ws.Range("O18").Value = 25.15362115366423

# Row 19
$ws.Range("C19").Value = 5.447322814145106
$ws.Range("D19").Value = 3.364975394290782
$ws.Range("E19").Value = 19.6212491263944
$ws.Range("F19").Value = 29.80695963907844
$ws.Range("G19").Value = 3.55531674854207
$ws.Range("N19").Value = 20.73962067985785
$ws.Range("O19").Value = 25.09721205081313

# Row 20
$ws.Range("C20").Value = 5.51561580698003
$ws.Range("D20").Value = 3.36948034857057
$ws.Range("E20").Value = 20.0030291576616
$ws.Range("F20").Value = 30.1528966668276
$ws.Range("G20").Value = 3.553733553664874
$ws.Range("N20").Value = 20.96544799483449
$ws.Range("O20").Value = 25.35062412894835

# Row 21
$ws.Range("C21").Value = 5.741001641805984
$ws.Range("D21").Value = 3.385785728614525
$ws.Range("E21").Value = 21.23863512412181
$ws.Range("F21").Value = 31.29822367207161
$ws.Range("G21").Value = 3.548570817188154
$ws.Range("N21").Value = 21.70751365554064
$ws.Range("O21").Value = 26.19358587191

# Row 22
$ws.Range("C22").Value = 5.885037196549976
$ws.Range("D22").Value = 3.397312611908086
$ws.Range("E22").Value = 22.01125871943847
$ws.Range("F22").Value = 32.03284516023606
$ws.Range("G22").Value = 3.54531255446293
$ws.Range("N22").Value = 22.17935961385674
$ws.Range("O22").Value = 26.73721587234358

# Row 23
$ws.Range("C23").Value = 5.808475767081984
$ws.Range("D23").Value = 3.391080969226284
$ws.Range("E23").Value = 21.6020766720854
$ws.Range("F23").Value = 31.64211395879047
$ws.Range("G23").Value = 3.547040999811102
$ws.Range("N23").Value = 21.92877110911181
$ws.Range("O23").Value = 26.44779810753272

# Row 24
$ws.Range("C24").Value = 5.511880155143001
$ws.Range("D24").Value = 3.369228557453476
$ws.Range("E24").Value = 19.98224142129517
$ws.Range("F24").Value = 30.13395990796412
$ws.Range("G24").Value = 3.553819904017697
$ws.Range("N24").Value = 20.95310750188673
$ws.Range("O24").Value = 25.33673718956975

# Row 25
$ws.Range("C25").Value = 5.179520665868902
$ws.Range("D25").Value = 3.349398448047916
$ws.Range("E25").Value = 18.08374502840279
$ws.Range("F25").Value = 28.45597794450956
$ws.Range("G25").Value = 3.561634903240105
$ws.Range("N25").Value = 19.84905939529497
$ws.Range("O25").Value = 24.11354218685652
